# NYPD CompStat weekly sheet refresh: bump the report header (volume/number
# and the covered week's date range), then swap in the newly collected
# per-precinct crime figures for the week.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# Header text: "Volume 32   Number  14" -> "...Number  15"
# ---------------------------------------------------------------------
$ws.Range("A8").Characters(21, 2).Text = "15"

# ---------------------------------------------------------------------
# Header text: "Report Covering the Week  3/31/2025  Through  4/6/2025"
#           -> "Report Covering the Week  4/7/2025  Through  4/13/2025"
# Replace the later substring first so the earlier substring's start
# offset isn't shifted by a length change.
# ---------------------------------------------------------------------
$ws.Range("C9").Characters(47, 8).Text = "4/13/2025"
$ws.Range("C9").Characters(27, 9).Text = "4/7/2025"

# ---------------------------------------------------------------------
# Helper: force a cell to hold literal text (not an auto-coerced number)
# while keeping the clean, pre-existing cell style (no quote-prefix xf).
# ---------------------------------------------------------------------
function Set-TextValue($ws, $targetRef, $text, $styleSourceRef) {
    $ws.Range($targetRef).Value = "'" + $text
    $ws.Range($styleSourceRef).Copy() | Out-Null
    $ws.Range($targetRef).PasteSpecial(-4122) | Out-Null
}

# Helper: force a cell to hold a real number while reusing an existing
# numeric cell's clean style.
function Set-NumberValue($ws, $targetRef, $number, $styleSourceRef) {
    $ws.Range($targetRef).Value = $number
    $ws.Range($styleSourceRef).Copy() | Out-Null
    $ws.Range($targetRef).PasteSpecial(-4122) | Out-Null
}

# ---------------------------------------------------------------------
# Row 15 (Rape): D/E flip from numbers to "0" / "***.*" text markers.
# ---------------------------------------------------------------------
Set-TextValue $ws "D15" "0" "C14"
Set-TextValue $ws "E15" "***.*" "E14"
$ws.Range("F15").Value = 1
$ws.Range("G15").Value = 3

# Row 16 (Robbery)
$ws.Range("C16").Value = 7
$ws.Range("D16").Value = 8
$ws.Range("E16").Value = -12.5
$ws.Range("I16").Value = 87
$ws.Range("J16").Value = 95
$ws.Range("K16").Value = -8.421052631578
$ws.Range("L16").Value = -40
$ws.Range("M16").Value = 74
$ws.Range("N16").Value = -87.482014388489

# Row 17 (Fel. Assault)
$ws.Range("C17").Value = 6
$ws.Range("D17").Value = 8
$ws.Range("E17").Value = -25
$ws.Range("F17").Value = 43
$ws.Range("H17").Value = 4.878048780487
$ws.Range("I17").Value = 147
$ws.Range("J17").Value = 159
$ws.Range("K17").Value = -7.547169811320
$ws.Range("L17").Value = -8.125
$ws.Range("M17").Value = 149.152542372881
$ws.Range("N17").Value = -22.631578947368

# Row 18 (Burglary)
$ws.Range("C18").Value = 8
$ws.Range("D18").Value = 2
$ws.Range("E18").Value = 300
$ws.Range("F18").Value = 33
$ws.Range("G18").Value = 22
$ws.Range("H18").Value = 50
$ws.Range("I18").Value = 129
$ws.Range("J18").Value = 109
$ws.Range("K18").Value = 18.348623853211
$ws.Range("L18").Value = -1.526717557251
$ws.Range("M18").Value = 30.303030303030
$ws.Range("N18").Value = -83.793969849246

# Row 19 (Gr. Larceny)
$ws.Range("C19").Value = 21
$ws.Range("D19").Value = 36
$ws.Range("E19").Value = -41.666666666666
$ws.Range("F19").Value = 132
$ws.Range("G19").Value = 131
$ws.Range("H19").Value = 0.763358778625
$ws.Range("I19").Value = 508
$ws.Range("J19").Value = 611
$ws.Range("K19").Value = -16.857610474631
$ws.Range("L19").Value = -25.073746312684
$ws.Range("M19").Value = -19.108280254777
$ws.Range("N19").Value = -81.507098653076

# Row 20 (G.L.A.)
$ws.Range("J20").Value = 13
$ws.Range("K20").Value = -53.846153846153
$ws.Range("L20").Value = -60
$ws.Range("N20").Value = -94.594594594594

# Row 21 (TOTAL)
$ws.Range("C21").Value = 42
$ws.Range("D21").Value = 55
$ws.Range("E21").Value = -23.636363636363
$ws.Range("F21").Value = 234
$ws.Range("G21").Value = 229
$ws.Range("H21").Value = 2.183406113537
$ws.Range("I21").Value = 888
$ws.Range("J21").Value = 1000
$ws.Range("K21").Value = -11.2
$ws.Range("L21").Value = -21.554770318021
$ws.Range("M21").Value = 5.213270142180
$ws.Range("N21").Value = -80.479226203561

# Row 22 (Transit)
$ws.Range("C22").Value = 5
$ws.Range("D22").Value = 5
$ws.Range("E22").Value = 0
$ws.Range("G22").Value = 17
$ws.Range("H22").Value = -5.882352941176
$ws.Range("I22").Value = 52
$ws.Range("J22").Value = 55
$ws.Range("K22").Value = -5.454545454545
$ws.Range("L22").Value = -20
$ws.Range("M22").Value = 30

# Row 24 (Petit Larceny)
$ws.Range("C24").Value = 65
$ws.Range("D24").Value = 80
$ws.Range("E24").Value = -18.75
$ws.Range("F24").Value = 332
$ws.Range("G24").Value = 311
$ws.Range("H24").Value = 6.752411575562
$ws.Range("I24").Value = 1137
$ws.Range("J24").Value = 1172
$ws.Range("K24").Value = -2.986348122866
$ws.Range("L24").Value = 13.134328358209
$ws.Range("M24").Value = -10.823529411764

# Row 25 (Retail Theft)
$ws.Range("C25").Value = 59
$ws.Range("D25").Value = 74
$ws.Range("E25").Value = -20.270270270270
$ws.Range("F25").Value = 296
$ws.Range("G25").Value = 291
$ws.Range("H25").Value = 1.718213058419
$ws.Range("I25").Value = 1014
$ws.Range("J25").Value = 1025
$ws.Range("K25").Value = -1.073170731707
$ws.Range("L25").Value = 8.102345415778

# Row 26 (Misd. Assault)
$ws.Range("C26").Value = 18
$ws.Range("D26").Value = 20
$ws.Range("E26").Value = -10
$ws.Range("F26").Value = 88
$ws.Range("G26").Value = 67
$ws.Range("H26").Value = 31.343283582089
$ws.Range("I26").Value = 308
$ws.Range("J26").Value = 284
$ws.Range("K26").Value = 8.450704225352
$ws.Range("L26").Value = 13.653136531365
$ws.Range("M26").Value = 62.962962962963

# Row 27 (UCR Rape*): D/E flip from numbers to "0" / "***.*" text markers.
Set-TextValue $ws "D27" "0" "C14"
Set-TextValue $ws "E27" "***.*" "E14"
$ws.Range("F27").Value = 1
$ws.Range("G27").Value = 4
$ws.Range("H27").Value = -75

# Row 28 (Other Sex Crimes)
$ws.Range("C28").Value = 3
$ws.Range("D28").Value = 8
$ws.Range("E28").Value = -62.5
$ws.Range("F28").Value = 12
$ws.Range("G28").Value = 24
$ws.Range("H28").Value = -50
$ws.Range("I28").Value = 57
$ws.Range("J28").Value = 60
$ws.Range("K28").Value = -5
$ws.Range("L28").Value = -12.307692307692

# Row 31 (Hate Crimes): D/E flip from "0" / "***.*" text markers to numbers.
Set-NumberValue $ws "D31" 1 "J14"
Set-NumberValue $ws "E31" -100 "K14"
$ws.Range("G31").Value = 2
$ws.Range("I31").Value = 4
$ws.Range("J31").Value = 6
$ws.Range("K31").Value = -33.333333333333
$ws.Range("L31").Value = -42.857142857142

# Row 33 (Traffic Fatalities): C flips from a number to the "0" text marker.
Set-TextValue $ws "C33" "0" "C14"

$excel.CutCopyMode = $false
